# "update sd card settings"
#
# Updates the "tanks" sheet: replaces the placeholder tank_max_volume
# (column C) values with the real, measured volumes, and switches the
# sensor_mounting_height column (F) from the text placeholder "10" to a
# real numeric reading (15 for the one outlier row, 10 for the rest).
# Also reproduces the view-state differences recorded for the sheets
# (active sheet/selection, a handful of resized columns on "tanks").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("valves")
$ws2 = $wb.Worksheets.Item("tanks")

# --- tanks!C2:C29 (tank_max_volume) and tanks!F2:F29 (sensor_mounting_height) ---

$ws2.Cells.Item(2,3).Value  = 1307.7
$ws2.Cells.Item(2,6).Value  = 10

$ws2.Cells.Item(3,3).Value  = 1307.7
$ws2.Cells.Item(3,6).Value  = 10

$ws2.Cells.Item(4,3).Value  = 1737.6
$ws2.Cells.Item(4,6).Value  = 10

$ws2.Cells.Item(5,3).Value  = 1737.6
$ws2.Cells.Item(5,6).Value  = 10

$ws2.Cells.Item(6,3).Value  = 3169.1
$ws2.Cells.Item(6,6).Value  = 10

$ws2.Cells.Item(7,3).Value  = 3169.1
$ws2.Cells.Item(7,6).Value  = 10

$ws2.Cells.Item(8,3).Value  = 3149.6
$ws2.Cells.Item(8,6).Value  = 10

$ws2.Cells.Item(9,3).Value  = 3222.9
$ws2.Cells.Item(9,6).Value  = 10

$ws2.Cells.Item(10,3).Value = 839.5
$ws2.Cells.Item(10,6).Value = 10

$ws2.Cells.Item(11,3).Value = 839.5
$ws2.Cells.Item(11,6).Value = 10

$ws2.Cells.Item(12,3).Value = 1508.4
$ws2.Cells.Item(12,6).Value = 10

$ws2.Cells.Item(13,3).Value = 661.7
$ws2.Cells.Item(13,6).Value = 10

$ws2.Cells.Item(14,3).Value = 396.5
$ws2.Cells.Item(14,6).Value = 10

$ws2.Cells.Item(15,3).Value = 631.29999999999995
$ws2.Cells.Item(15,6).Value = 10

$ws2.Cells.Item(16,3).Value = 308.2
$ws2.Cells.Item(16,6).Value = 10

$ws2.Cells.Item(17,3).Value = 367.3
$ws2.Cells.Item(17,6).Value = 10

$ws2.Cells.Item(18,3).Value = 264.39999999999998
$ws2.Cells.Item(18,6).Value = 10

$ws2.Cells.Item(19,3).Value = 24.2
$ws2.Cells.Item(19,6).Value = 10

$ws2.Cells.Item(20,3).Value = 24.2
$ws2.Cells.Item(20,6).Value = 10

$ws2.Cells.Item(21,3).Value = 18.2
$ws2.Cells.Item(21,6).Value = 10

$ws2.Cells.Item(22,3).Value = 18.2
$ws2.Cells.Item(22,6).Value = 10

$ws2.Cells.Item(23,3).Value = 24.8
$ws2.Cells.Item(23,6).Value = 10

$ws2.Cells.Item(24,3).Value = 35.799999999999997
$ws2.Cells.Item(24,6).Value = 10

$ws2.Cells.Item(25,3).Value = 483.6
$ws2.Cells.Item(25,6).Value = 10

$ws2.Cells.Item(26,3).Value = 36
$ws2.Cells.Item(26,6).Value = 10

$ws2.Cells.Item(27,3).Value = 68.5
$ws2.Cells.Item(27,6).Value = 15

$ws2.Cells.Item(28,3).Value = 209.7
$ws2.Cells.Item(28,6).Value = 10

$ws2.Cells.Item(29,3).Value = 209.7
$ws2.Cells.Item(29,6).Value = 10

# --- column width tweaks on "tanks" (C, G, H, K, L, M got resized) ---

$ws2.Columns.Item(3).ColumnWidth  = 20.833333333333332
$ws2.Columns.Item(7).ColumnWidth  = 13.166666666666666
$ws2.Columns.Item(8).ColumnWidth  = 12.0
$ws2.Columns.Item(11).ColumnWidth = 12.5
$ws2.Columns.Item(12).ColumnWidth = 12.0
$ws2.Columns.Item(13).ColumnWidth = 11.166666666666666

# --- view state: "tanks" selection moves to G34 ... ---
$ws2.Range("G34").Select()

# --- ... but "valves" becomes the active/selected sheet, selection F54 ---
$ws1.Activate()
$ws1.Range("F54").Select()
